# Add a new "Yearly demand" worksheet at the end of the workbook.
#
# The sheet mirrors the layout already used by sheets like "DG Dispatch" /
# "Connected Households": a header row (B1:Y1 = 0..23) and an index
# column A (A2=0, A3=1, A4=2) for the three data rows, all styled with
# the existing bold/bordered/centered header style.
#
# Easiest (and most faithful) way to get that skeleton - including its
# styles, sheetPr and page margins - is to duplicate an existing sheet
# that already has it ("DG Dispatch") and then overwrite its data.

$wb = $excel.ActiveWorkbook

$template = $wb.Worksheets.Item("DG Dispatch")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$template.Copy($null, $lastSheet)

# The copy lands right after the last existing sheet and becomes active;
# grab it by its default copy name and rename it.
$ws = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws.Name = "Yearly demand"

# Restore the original active sheet/selection.
$template.Activate()

# --- Overwrite the demand data (header row + index column already match
# the template, so only the B:Y data rows need to be replaced). ---------
$row2 = @(-32.5,-19.5,-13,-13,-13,142.5,291.5,327,388.5,502,596,670.5,745,651,576.5,502,320.5,139,32,-117,-97.5,-78,-52,-39)
$row3 = @(-32.5,-19.5,-13,0,0,-19.5,0,324,486,648,729,751.5,583,567,333.5,340,243,57.99999999999999,-130,0,0,-78,0,-39)
$row4 = @(-32.5,-19.5,0,0,0,-19.5,0,0,81,324,567,589.5,648,567,324,162,81,0,-130,0,0,0,0,-39)

for ($i = 0; $i -lt $row2.Length; $i++) {
    $ws.Cells.Item(2, $i + 2).Value = $row2[$i]
}
for ($i = 0; $i -lt $row3.Length; $i++) {
    $ws.Cells.Item(3, $i + 2).Value = $row3[$i]
}
for ($i = 0; $i -lt $row4.Length; $i++) {
    $ws.Cells.Item(4, $i + 2).Value = $row4[$i]
}
